$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): column C becomes old D's "prediction" text,
# column D becomes old E's "rejection-f" text, column E removed entirely.
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"

# Row 2: new B2 value, column C takes previous D2's text ("g__CAG-631"),
# D2 keeps "g__CAG-631" (previously E2), E2 removed.
$ws.Range("B2").Value = 0.0871431380098926
$ws.Range("C2").Value = "g__CAG-631"
$ws.Range("D2").Value = "g__CAG-631"

# Row 3: new B3 value, column C takes "g__CAG-631", D3 becomes "g__CAG-631(reject)".
$ws.Range("B3").Value = -0.04862773867910164
$ws.Range("C3").Value = "g__CAG-631"
$ws.Range("D3").Value = "g__CAG-631(reject)"

# Remove the now-unused column E entirely so the sheet dimension shrinks to A1:D3.
$ws.Range("E1:E3").Delete()
